$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Förändrad) from row 2 to row 296: 45181 -> 45182
$ws.Range("C2:C296").Value = 45182
